# Updated cryptos list (Price + Volume(1h) columns) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without letting Excel
# auto-coerce numeric-looking strings (e.g. "330.66") into floating point
# numbers. We briefly force a text NumberFormat, assign the value, then
# restore the "Normal" style so the cell format matches the source file
# (plain text, default/unstyled cell) exactly.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "30.459.71"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.092.92"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "330.66"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  +0.16%  "
Set-TextValue "D7" "0.5214"
$ws.Range("E7").Value = "  -0.47%  "
Set-TextValue "D8" "0.4383"
$ws.Range("E8").Value = "  -0.62%  "
Set-TextValue "D9" "54.12"
$ws.Range("E9").Value = "  +15.45%  "
Set-TextValue "D10" "0.08914"
$ws.Range("E10").Value = "  -1.70%  "
Set-TextValue "D11" "1.154"
$ws.Range("E11").Value = "  -2.72%  "
Set-TextValue "D12" "24.28"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").Value = "2.086.47"
$ws.Range("E13").Value = "  -1.54%  "
Set-TextValue "D14" "6.693"
$ws.Range("E14").Value = "  -1.28%  "
Set-TextValue "D15" "7.681"
$ws.Range("E15").Value = "  -2.46%  "
Set-TextValue "D16" "95.93"
$ws.Range("E16").Value = "  -2.43%  "
Set-TextValue "D17" "1.003"
$ws.Range("E17").Value = "  +0.16%  "
Set-TextValue "D18" "0.00001122"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("E19").Value = "  -0.53%  "
Set-TextValue "D20" "19.17"
$ws.Range("E20").Value = "  -0.27%  "
Set-TextValue "D21" "1.001"
$ws.Range("E21").Value = "  +0.15%  "
Set-TextValue "D22" "6.263"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "30.487.40"
$ws.Range("E24").Value = "  +1.22%  "
Set-TextValue "D25" "2.324"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").Value = "2.329.66"
$ws.Range("E26").Value = "  -1.63%  "
Set-TextValue "D27" "22.26"
$ws.Range("E27").Value = "  -3.14%  "
Set-TextValue "D28" "2.556"
$ws.Range("E28").Value = "  -0.45%  "
Set-TextValue "D29" "163.13"
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue "D30" "131.62"
$ws.Range("E30").Value = "  -1.43%  "
Set-TextValue "D31" "1.183"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -0.25%  "
Set-TextValue "D33" "1.658"
$ws.Range("E33").Value = "  +7.06%  "
Set-TextValue "D34" "6.160"
$ws.Range("E34").Value = "  -1.59%  "
Set-TextValue "D35" "3.900"
$ws.Range("E35").Value = "  -3.16%  "
Set-TextValue "D36" "10.03"
$ws.Range("E36").Value = "  +4.79%  "
Set-TextValue "D37" "0.02563"
$ws.Range("E37").Value = "  -1.50%  "
Set-TextValue "D38" "0.06831"
$ws.Range("E38").Value = "  +0.92%  "
Set-TextValue "D39" "5.476"
$ws.Range("E39").Value = "  -2.15%  "
Set-TextValue "D40" "12.60"
$ws.Range("E40").Value = "  -1.26%  "
Set-TextValue "D41" "0.2255"
$ws.Range("E41").Value = "  -0.94%  "
Set-TextValue "D42" "0.6878"
$ws.Range("E42").Value = "  +0.46%  "
Set-TextValue "D43" "1.253"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  +0.15%  "
Set-TextValue "D45" "14.01"
$ws.Range("E45").Value = "  -0.86%  "
Set-TextValue "D46" "0.6340"
$ws.Range("E46").Value = "  -1.44%  "
Set-TextValue "D47" "2.195"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("E48").Value = "  -1.39%  "
Set-TextValue "D49" "1.233"
$ws.Range("E49").Value = "  +6.53%  "
Set-TextValue "D50" "1.243"
$ws.Range("E50").Value = "  -3.41%  "
Set-TextValue "D51" "81.80"
$ws.Range("E51").Value = "  -1.72%  "
